# "ya guarda todo con estados" -- add a new "estado" (status) column (J)
# to the sheet: header in J1, "disponible" for the first listing (row 2),
# and an (empty but styled, like the rest of row 3) cell for the second
# listing (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header + value
$ws.Range("J1").Value = "estado"
$ws.Range("J2").Value = "disponible"
# J3 is intentionally left blank (no estado recorded yet for that listing)
# but still needs the row's formatting, handled by the PasteSpecial below.

# Match the look & feel of each row: copy the formatting only (no values)
# from an existing cell in the same row so the new column blends in with
# the rest of the table instead of using Excel's plain default style.
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("H2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("H3").Copy() | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# PasteSpecial only carries formats, so re-assert the text values that
# matter (and make sure J3 truly stays empty).
$ws.Range("J1").Value = "estado"
$ws.Range("J2").Value = "disponible"
$ws.Range("J3").Value = $null

# Scroll the window right so column F becomes the left-most visible
# column, then land the selection on J4 (just below the new column),
# mirroring how the sheet was left after the edit.
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("J4").Select()
